$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 928.8570999999999
$ws.Range("J121").Value = 1188.0667
$ws.Range("L121").Value = 3564.2001
$ws.Range("N121").Value = -7058.2001
$ws.Range("H132").Value = 12145.87
$ws.Range("I132").Value = 1296.3472
$ws.Range("J132").Value = 72235.53999999999
$ws.Range("K132").Value = 3889.0416
$ws.Range("L132").Value = 216706.62
$ws.Range("M132").Value = -1359.0416
$ws.Range("N132").Value = -221766.62
$ws.Range("H137").Value = 2299.6626
$ws.Range("I137").Value = 693.5484
$ws.Range("J137").Value = 7041.524
$ws.Range("K137").Value = 2080.6452
$ws.Range("L137").Value = 21124.572
$ws.Range("M137").Value = 469.3548000000001
$ws.Range("N137").Value = -26224.572
$ws.Range("H138").Value = 3660.0986
$ws.Range("I138").Value = 2225.125
$ws.Range("J138").Value = 4392.851
$ws.Range("K138").Value = 6675.375
$ws.Range("L138").Value = 13178.553
$ws.Range("M138").Value = -1535.375
$ws.Range("N138").Value = -23458.553

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1768.84
$ws.Range("I32").Value = 1214.3146
$ws.Range("J32").Value = 6255.4546
$ws.Range("K32").Value = 1214.3146
$ws.Range("L32").Value = 6255.4546
$ws.Range("M32").Value = -927.3145999999999
$ws.Range("N32").Value = -6829.4546
$ws.Range("H61").Value = 1141.8062
$ws.Range("I61").Value = 970.95776
$ws.Range("J61").Value = 1591.0741
$ws.Range("K61").Value = 970.95776
$ws.Range("L61").Value = 1591.0741
$ws.Range("M61").Value = -758.95776
$ws.Range("N61").Value = -2015.0741
$ws.Range("H74").Value = 1794.1954
$ws.Range("I74").Value = 1661.24
$ws.Range("J74").Value = 2625.1667
$ws.Range("K74").Value = 1661.24
$ws.Range("L74").Value = 2625.1667
$ws.Range("M74").Value = -787.24
$ws.Range("N74").Value = -4373.1667
$ws.Range("H77").Value = 1794.1954
$ws.Range("I77").Value = 1661.24
$ws.Range("J77").Value = 2625.1667
$ws.Range("K77").Value = 8306.200000000001
$ws.Range("L77").Value = 13125.8335
$ws.Range("M77").Value = -3938.200000000001
$ws.Range("N77").Value = -21861.8335
$ws.Range("H97").Value = 531.8182
$ws.Range("I97").Value = 418.02777
$ws.Range("K97").Value = 418.02777
$ws.Range("M97").Value = 77.97223000000002
$ws.Range("H132").Value = 6668149
$ws.Range("I132").Value = 10205094
$ws.Range("J132").Value = 2368.8845
$ws.Range("K132").Value = 30615282
$ws.Range("L132").Value = 7106.6535
$ws.Range("M132").Value = -30612752
$ws.Range("N132").Value = -12166.6535
$ws.Range("H136").Value = 1141.8062
$ws.Range("I136").Value = 970.95776
$ws.Range("J136").Value = 1591.0741
$ws.Range("K136").Value = 2912.87328
$ws.Range("L136").Value = 4773.2223
$ws.Range("M136").Value = -362.8732799999998
$ws.Range("N136").Value = -9873.222300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2217.8333
$ws.Range("I105").Value = 2570
$ws.Range("J105").Value = 2147.4
$ws.Range("K105").Value = 2570
$ws.Range("L105").Value = 2147.4
$ws.Range("M105").Value = -823
$ws.Range("N105").Value = -5641.4
$ws.Range("H134").Value = 2376.18
$ws.Range("I134").Value = 916.86957
$ws.Range("J134").Value = 3619.2964
$ws.Range("K134").Value = 2750.60871
$ws.Range("L134").Value = 10857.8892
$ws.Range("M134").Value = -215.60871
$ws.Range("N134").Value = -15927.8892

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 16667450
$ws.Range("I7").Value = 20000740
$ws.Range("K7").Value = 20000740
$ws.Range("M7").Value = -20000627
$ws.Range("H17").Value = 28375
$ws.Range("I17").Value = 27750
$ws.Range("J17").Value = 29000
$ws.Range("K17").Value = 27750
$ws.Range("L17").Value = 29000
$ws.Range("M17").Value = -27576
$ws.Range("N17").Value = -29348
$ws.Range("H31").Value = 2075.62
$ws.Range("I31").Value = 891.8
$ws.Range("J31").Value = 3522.5112
$ws.Range("K31").Value = 891.8
$ws.Range("L31").Value = 3522.5112
$ws.Range("M31").Value = -596.8
$ws.Range("N31").Value = -4112.5112
$ws.Range("H34").Value = 2075.62
$ws.Range("I34").Value = 891.8
$ws.Range("J34").Value = 3522.5112
$ws.Range("K34").Value = 891.8
$ws.Range("L34").Value = 3522.5112
$ws.Range("M34").Value = -689.8
$ws.Range("N34").Value = -3926.5112
$ws.Range("H58").Value = 1191.3959
$ws.Range("I58").Value = 794.8
$ws.Range("J58").Value = 3174.375
$ws.Range("K58").Value = 794.8
$ws.Range("L58").Value = 3174.375
$ws.Range("M58").Value = -591.8
$ws.Range("N58").Value = -3580.375
$ws.Range("H99").Value = 1859.7727
$ws.Range("J99").Value = 1849.8823
$ws.Range("L99").Value = 1849.8823
$ws.Range("N99").Value = -4845.8823
$ws.Range("H126").Value = 1859.7727
$ws.Range("J126").Value = 1849.8823
$ws.Range("L126").Value = 5549.6469
$ws.Range("N126").Value = -10489.6469
$ws.Range("H132").Value = 22946.445
$ws.Range("I132").Value = 1006.1071
$ws.Range("J132").Value = 159464.11
$ws.Range("K132").Value = 3018.3213
$ws.Range("L132").Value = 478392.33
$ws.Range("M132").Value = -488.3212999999996
$ws.Range("N132").Value = -483452.33
$ws.Range("H134").Value = 246932.8
$ws.Range("I134").Value = 814.32556
$ws.Range("J134").Value = 1002868.06
$ws.Range("K134").Value = 2442.97668
$ws.Range("L134").Value = 3008604.18
$ws.Range("M134").Value = 92.02332000000024
$ws.Range("N134").Value = -3013674.18
$ws.Range("H136").Value = 1191.3959
$ws.Range("I136").Value = 794.8
$ws.Range("J136").Value = 3174.375
$ws.Range("K136").Value = 2384.4
$ws.Range("L136").Value = 9523.125
$ws.Range("M136").Value = 165.6000000000004
$ws.Range("N136").Value = -14623.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 3507
$ws.Range("I87").Value = 3507
$ws.Range("K87").Value = 10521
$ws.Range("M87").Value = -9273
$ws.Range("H90").Value = 3507
$ws.Range("I90").Value = 3507
$ws.Range("K90").Value = 31563
$ws.Range("M90").Value = -25323
$ws.Range("H131").Value = 2080.6948
$ws.Range("I131").Value = 10440.9
$ws.Range("J131").Value = 1097.1412
$ws.Range("K131").Value = 31322.7
$ws.Range("L131").Value = 3291.4236
$ws.Range("M131").Value = -26282.7
$ws.Range("N131").Value = -13371.4236
$ws.Range("H132").Value = 3582.375
$ws.Range("I132").Value = 1175
$ws.Range("J132").Value = 5989.75
$ws.Range("K132").Value = 10575
$ws.Range("L132").Value = 53907.75
$ws.Range("M132").Value = -8045
$ws.Range("N132").Value = -58967.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12417.667
$ws.Range("I43").Value = 5341.3335
$ws.Range("J43").Value = 13428.571
$ws.Range("K43").Value = 5341.3335
$ws.Range("L43").Value = 13428.571
$ws.Range("M43").Value = -5190.3335
$ws.Range("N43").Value = -13730.571
$ws.Range("H46").Value = 22943.8
$ws.Range("J46").Value = 22943.8
$ws.Range("L46").Value = 22943.8
$ws.Range("N46").Value = -23255.8
$ws.Range("H80").Value = 203687.4
$ws.Range("I80").Value = 360978.94
$ws.Range("J80").Value = 3498.182
$ws.Range("K80").Value = 360978.94
$ws.Range("L80").Value = 3498.182
$ws.Range("M80").Value = -359980.94
$ws.Range("N80").Value = -5494.182
$ws.Range("H83").Value = 203687.4
$ws.Range("I83").Value = 360978.94
$ws.Range("J83").Value = 3498.182
$ws.Range("K83").Value = 1804894.7
$ws.Range("L83").Value = 17490.91
$ws.Range("M83").Value = -1799902.7
$ws.Range("N83").Value = -27474.91
$ws.Range("H96").Value = 30993.334
$ws.Range("J96").Value = 30993.334
$ws.Range("L96").Value = 30993.334
$ws.Range("N96").Value = -36485.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 992.125
$ws.Range("I22").Value = 937.5
$ws.Range("J22").Value = 1046.75
$ws.Range("K22").Value = 937.5
$ws.Range("L22").Value = 1046.75
$ws.Range("M22").Value = -642.5
$ws.Range("N22").Value = -1636.75
$ws.Range("H27").Value = 992.125
$ws.Range("I27").Value = 937.5
$ws.Range("J27").Value = 1046.75
$ws.Range("K27").Value = 937.5
$ws.Range("L27").Value = 1046.75
$ws.Range("M27").Value = -830.5
$ws.Range("N27").Value = -1260.75
$ws.Range("H55").Value = 395.05
$ws.Range("I55").Value = 271.9375
$ws.Range("J55").Value = 887.5
$ws.Range("K55").Value = 271.9375
$ws.Range("L55").Value = 887.5
$ws.Range("M55").Value = -98.9375
$ws.Range("N55").Value = -1233.5
$ws.Range("H132").Value = 2171.7544
$ws.Range("I132").Value = 1439.6945
$ws.Range("J132").Value = 3426.7144
$ws.Range("K132").Value = 4319.083500000001
$ws.Range("L132").Value = 10280.1432
$ws.Range("M132").Value = -1789.083500000001
$ws.Range("N132").Value = -15340.1432
$ws.Range("H136").Value = 1507.4559
$ws.Range("I136").Value = 1147.22
$ws.Range("J136").Value = 2508.111
$ws.Range("K136").Value = 3441.66
$ws.Range("L136").Value = 7524.333
$ws.Range("M136").Value = -891.6599999999999
$ws.Range("N136").Value = -12624.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 805.31146
$ws.Range("I136").Value = 637.0333000000001
$ws.Range("J136").Value = 968.1613
$ws.Range("K136").Value = 1911.0999
$ws.Range("L136").Value = 2904.4839
$ws.Range("M136").Value = 638.9000999999998
$ws.Range("N136").Value = -8004.4839
